# Auto-generated Excel COM-interop script applying scheduled-runner market price updates
# across all 8 leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 87705
$ws.Range("I43").Value = 919.25
$ws.Range("J43").Value = 131097.88
$ws.Range("K43").Value = 919.25
$ws.Range("L43").Value = 131097.88
$ws.Range("M43").Value = -850.25
$ws.Range("N43").Value = -131235.88
$ws.Range("H110").Value = 64220
$ws.Range("J110").Value = 64220
$ws.Range("L110").Value = 64220
$ws.Range("N110").Value = -72400
$ws.Range("H112").Value = 346462.12
$ws.Range("J112").Value = 358803.66
$ws.Range("L112").Value = 1076410.98
$ws.Range("N112").Value = -1078626.98
$ws.Range("H120").Value = 49993.332
$ws.Range("J120").Value = 49993.332
$ws.Range("L120").Value = 49993.332
$ws.Range("N120").Value = -59669.332
$ws.Range("H123").Value = 88922.5
$ws.Range("J123").Value = 88922.5
$ws.Range("L123").Value = 88922.5
$ws.Range("N123").Value = -98722.5
$ws.Range("H134").Value = 61791.363
$ws.Range("J134").Value = 65881.11
$ws.Range("L134").Value = 65881.11
$ws.Range("N134").Value = -76021.11
$ws.Range("H135").Value = 833
$ws.Range("I135").Value = 492.89474
$ws.Range("K135").Value = 4436.05266
$ws.Range("M135").Value = -1901.05266
$ws.Range("H136").Value = 72867.7
$ws.Range("J136").Value = 82084.625
$ws.Range("L136").Value = 82084.625
$ws.Range("N136").Value = -92284.625
$ws.Range("H138").Value = 2270.8157
$ws.Range("I138").Value = 1562.7142
$ws.Range("K138").Value = 4688.142599999999
$ws.Range("M138").Value = 451.8574000000008

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6790.7144
$ws.Range("I45").Value = 7089.222
$ws.Range("K45").Value = 7089.222
$ws.Range("M45").Value = -6712.222
$ws.Range("H52").Value = 49785.332
$ws.Range("J52").Value = 49785.332
$ws.Range("L52").Value = 49785.332
$ws.Range("N52").Value = -50421.332
$ws.Range("H108").Value = 89992
$ws.Range("J108").Value = 89992
$ws.Range("L108").Value = 89992
$ws.Range("N108").Value = -97672
$ws.Range("H115").Value = 55524.4
$ws.Range("J115").Value = 66655.5
$ws.Range("L115").Value = 66655.5
$ws.Range("N115").Value = -69789.5
$ws.Range("H117").Value = 71100.75
$ws.Range("J117").Value = 71100.75
$ws.Range("L117").Value = 71100.75
$ws.Range("N117").Value = -80278.75
$ws.Range("H121").Value = 64609
$ws.Range("J121").Value = 64609
$ws.Range("L121").Value = 64609
$ws.Range("N121").Value = -68103
$ws.Range("H122").Value = 3015.3655
$ws.Range("I122").Value = 3045.7104
$ws.Range("K122").Value = 9137.1312
$ws.Range("M122").Value = -6687.1312

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 77678
$ws.Range("J51").Value = 77678
$ws.Range("L51").Value = 77678
$ws.Range("N51").Value = -78660
$ws.Range("H52").Value = 68422
$ws.Range("J52").Value = 68422
$ws.Range("L52").Value = 68422
$ws.Range("N52").Value = -68948
$ws.Range("H86").Value = 2558.0588
$ws.Range("I86").Value = 1669.8572
$ws.Range("K86").Value = 1669.8572
$ws.Range("M86").Value = -546.8571999999999
$ws.Range("H89").Value = 2558.0588
$ws.Range("I89").Value = 1669.8572
$ws.Range("K89").Value = 8349.286
$ws.Range("M89").Value = -2733.286
$ws.Range("H117").Value = 77054.75
$ws.Range("J117").Value = 77054.75
$ws.Range("L117").Value = 77054.75
$ws.Range("N117").Value = -86232.75
$ws.Range("H121").Value = 68422
$ws.Range("J121").Value = 68422
$ws.Range("L121").Value = 68422
$ws.Range("N121").Value = -71916
$ws.Range("H132").Value = 30666.111
$ws.Range("J132").Value = 30666.111
$ws.Range("L132").Value = 30666.111
$ws.Range("N132").Value = -40786.111
$ws.Range("H134").Value = 2173.95
$ws.Range("I134").Value = 1590.3636
$ws.Range("K134").Value = 4771.0908
$ws.Range("M134").Value = -2236.0908
$ws.Range("H135").Value = 101278.57
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("H138").Value = 99990
$ws.Range("J138").Value = 99990
$ws.Range("L138").Value = 99990
$ws.Range("N138").Value = -110270
$ws.Range("M135").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
$ws.Range("H108").Value = 42614.145
$ws.Range("J108").Value = 42614.145
$ws.Range("L108").Value = 42614.145
$ws.Range("N108").Value = -50294.145
$ws.Range("H116").Value = 40789.082
$ws.Range("J116").Value = 40789.082
$ws.Range("L116").Value = 40789.082
$ws.Range("N116").Value = -49967.082
$ws.Range("H119").Value = 96992.5
$ws.Range("J119").Value = 96992.5
$ws.Range("L119").Value = 96992.5
$ws.Range("N119").Value = -106668.5
$ws.Range("H132").Value = 1695.6765
$ws.Range("I132").Value = 1482.5518
$ws.Range("J132").Value = 2931.8
$ws.Range("K132").Value = 4447.6554
$ws.Range("L132").Value = 8795.400000000001
$ws.Range("M132").Value = -1917.6554
$ws.Range("N132").Value = -13855.4

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1619.7142
$ws.Range("J68").Value = 866.3333
$ws.Range("L68").Value = 2598.9999
$ws.Range("N68").Value = -4220.9999
$ws.Range("H71").Value = 1619.7142
$ws.Range("J71").Value = 866.3333
$ws.Range("L71").Value = 7796.9997
$ws.Range("N71").Value = -15908.9997
$ws.Range("H110").Value = 9500
$ws.Range("I110").Value = 8500
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 25500
$ws.Range("L110").Value = 30000
$ws.Range("M110").Value = -21410
$ws.Range("N110").Value = -38180
$ws.Range("H138").Value = 4754.6
$ws.Range("I138").Value = 2621.3333
$ws.Range("K138").Value = 7863.999899999999
$ws.Range("M138").Value = -2723.999899999999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 25148.834
$ws.Range("J52").Value = 25148.834
$ws.Range("L52").Value = 25148.834
$ws.Range("N52").Value = -25666.834
$ws.Range("H64").Value = 500045000
$ws.Range("I64").Value = 90000
$ws.Range("J64").Value = 1000000000
$ws.Range("K64").Value = 90000
$ws.Range("L64").Value = 1000000000
$ws.Range("M64").Value = -89752
$ws.Range("N64").Value = -1000000496
$ws.Range("H67").Value = 500045000
$ws.Range("I67").Value = 90000
$ws.Range("J67").Value = 1000000000
$ws.Range("K67").Value = 90000
$ws.Range("L67").Value = 1000000000
$ws.Range("M67").Value = -89142
$ws.Range("N67").Value = -1000001716
$ws.Range("H108").Value = 64940
$ws.Range("J108").Value = 64940
$ws.Range("L108").Value = 64940
$ws.Range("N108").Value = -72620
$ws.Range("H113").Value = 1729.5
$ws.Range("I113").Value = 1200
$ws.Range("K113").Value = 1200
$ws.Range("M113").Value = 970
$ws.Range("H122").Value = 9132.947
$ws.Range("I122").Value = 3619.1765
$ws.Range("J122").Value = 56000
$ws.Range("K122").Value = 10857.5295
$ws.Range("L122").Value = 168000
$ws.Range("M122").Value = -8407.529500000001
$ws.Range("N122").Value = -172900
$ws.Range("H132").Value = 2079.4285
$ws.Range("I132").Value = 1983
$ws.Range("J132").Value = 2236.125
$ws.Range("K132").Value = 5949
$ws.Range("L132").Value = 6708.375
$ws.Range("M132").Value = -3419
$ws.Range("N132").Value = -11768.375
$ws.Range("H135").Value = 69368
$ws.Range("J135").Value = 69368
$ws.Range("L135").Value = 69368
$ws.Range("N135").Value = -79508

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3000
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -798
$ws.Range("N61").Value = -5404
$ws.Range("H68").Value = 2902.111
$ws.Range("I68").Value = 2949.8333
$ws.Range("K68").Value = 2949.8333
$ws.Range("M68").Value = -2200.8333
$ws.Range("H71").Value = 2902.111
$ws.Range("I71").Value = 2949.8333
$ws.Range("K71").Value = 14749.1665
$ws.Range("M71").Value = -11005.1665
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -9340
$ws.Range("H136").Value = 5234.8184
$ws.Range("I136").Value = 5288.4
$ws.Range("K136").Value = 15865.2
$ws.Range("M136").Value = -13315.2

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 140136.75
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()

Write-Output "Applied scheduled market price updates across all sheets."